# Atualiza instrução de trabalho
# Rename the responsible person "Erick Silva" -> "Erick da Silva" on the
# ITI sheet (rows 3-24, column B / "Responsavel"), then leave the
# workbook with the ITI sheet active and cell E19 selected there (matching
# where the edit was made).

$wb = $excel.ActiveWorkbook

$iti = $wb.Worksheets.Item("ITI")

# Update the "Responsavel" column for the rows that referenced "Erick Silva".
$iti.Range("B3:B24").Value = "Erick da Silva"

# The edit was made on the ITI sheet, so it ends up the active tab/sheet,
# with E19 as the selected cell (matching the author's final cursor spot).
[void]$iti.Activate()
[void]$iti.Range("E19").Select()
